$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for team "The-Danish-German-Alliance" (old row 7) entirely
$ws.Rows.Item(7).Delete()

# Re-write the team / peer_group_1 / peer_group_2 assignment table to match the refreshed roster
$data = @(
    @('CH13FSQD', 'mms', 'monty-python-1'),
    @('EGPK', 'Numpy', 'hold-j1'),
    @('Inaugural-Project-2021_BSJH', 'ska', 'sft'),
    @('Numpy', 'mogens-pa-linkedin', 'lions'),
    @('a12', 'sev', 'm-m'),
    @('aristochats', 'team-kan-du-kode-eller-hvad', 'hello'),
    @('california-sunshine', 'ska', 'double-a'),
    @('concatenaters', 'kera', 'double-a'),
    @('cr-and-gb', 'california-sunshine', 'the-danish-german-alliance'),
    @('double-a', 'hold-j1', 'stronk'),
    @('econ-jedi', 'the-danish-german-alliance', 'egpk'),
    @('egpk', 'm-a-c-h', 'california-sunshine'),
    @('fgs', 'a12', 'mms'),
    @('group-sara-1', 'stronk', 'team-tons'),
    @('hello', 'concatenaters', 'monty-python-1'),
    @('hello-world', 'sstl', 'mf'),
    @('hold-i-nakken', 'hello-world', 'west-wing'),
    @('hold-j1', 'team-easy-on', 'sstl'),
    @('holdfast', 'unicorns_1', 'lions'),
    @('kera', 'Numpy', 'CH13FSQD'),
    @('lions', 'concatenaters', 'sft'),
    @('luna-hjerteknuser', 'sev', 'team-mej'),
    @('m-a-c-h', 'team-mej', 'regnefaetrene'),
    @('m-m', 'python-gurlz', 'CH13FSQD'),
    @('mf', 'python-gurlz', 'm-a-c-h'),
    @('mms', 'west-wing', 'hello-world'),
    @('mogens-pa-linkedin', 'fgs', 'unicorns_1'),
    @('monty-python-1', 'hold-i-nakken', 'team-olm'),
    @('programmorerne', 'a12', 'egpk'),
    @('python-gurlz', 'team-easy-on', 'hello'),
    @('regnefaetrene', 'hold-i-nakken', 'econ-jedi'),
    @('sev', 'team-olm', 'regnefaetrene'),
    @('sft', 'aristochats', 'team-tons'),
    @('ska', 'm-m', 'stressfaktoren'),
    @('sstl', 'kera', 'holdfast'),
    @('stressfaktoren', 'Inaugural-Project-2021_BSJH', 'econ-jedi'),
    @('stronk', 'mf', 'fgs'),
    @('team-easy-on', 'Inaugural-Project-2021_BSJH', 'programmorerne'),
    @('team-kan-du-kode-eller-hvad', 'luna-hjerteknuser', 'group-sara-1'),
    @('team-mej', 'team-kan-du-kode-eller-hvad', 'stressfaktoren'),
    @('team-olm', 'holdfast', 'group-sara-1'),
    @('team-tons', 'EGPK', 'mogens-pa-linkedin'),
    @('the-danish-german-alliance', 'luna-hjerteknuser', 'cr-and-gb'),
    @('unicorns_1', 'aristochats', 'programmorerne'),
    @('west-wing', 'EGPK', 'cr-and-gb'),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count